$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20240627-092035-"

# Bulk-update the "date" column (G) for every data row: 2024-06-26 -> 2024-06-27 (serial 45469 -> 45470)
$ws.Range("G2:G275").Value = 45470

# Row-specific balance corrections (columns D, E, H)
$ws.Cells.Item(17, 5).Value = 420.68
$ws.Cells.Item(17, 8).Value = 420.68
$ws.Cells.Item(97, 5).Value = 61278.51
$ws.Cells.Item(97, 8).Value = 61278.51
$ws.Cells.Item(101, 5).Value = 57053.42
$ws.Cells.Item(101, 8).Value = 57053.42
$ws.Cells.Item(104, 4).Value = -5229.71
$ws.Cells.Item(104, 5).Value = 5483.78
$ws.Cells.Item(104, 8).Value = 254.07
$ws.Cells.Item(112, 5).Value = 647.82000000000005
$ws.Cells.Item(112, 8).Value = 647.82000000000005
$ws.Cells.Item(113, 5).Value = 2093.59
$ws.Cells.Item(113, 8).Value = 2093.59
$ws.Cells.Item(118, 5).Value = 602.83000000000004
$ws.Cells.Item(118, 8).Value = 602.83000000000004
$ws.Cells.Item(138, 5).Value = 226.59
$ws.Cells.Item(138, 8).Value = 226.59
$ws.Cells.Item(158, 5).Value = 6097.68
$ws.Cells.Item(158, 8).Value = 6097.68
$ws.Cells.Item(171, 5).Value = 352.15
$ws.Cells.Item(171, 8).Value = 352.15
$ws.Cells.Item(172, 4).Value = 1531.91
$ws.Cells.Item(172, 8).Value = 10920.38
$ws.Cells.Item(173, 4).Value = -19625.599999999999
$ws.Cells.Item(173, 5).Value = 19899.28
$ws.Cells.Item(173, 8).Value = 273.68
$ws.Cells.Item(224, 5).Value = 396.23
$ws.Cells.Item(224, 8).Value = 396.23
$ws.Cells.Item(255, 5).Value = 985.28
$ws.Cells.Item(255, 8).Value = 985.28
$ws.Cells.Item(264, 4).Value = 0
$ws.Cells.Item(264, 5).Value = 733.08
$ws.Cells.Item(264, 8).Value = 733.08
$ws.Cells.Item(265, 4).Value = -14054.17
$ws.Cells.Item(265, 5).Value = 14317.24
$ws.Cells.Item(265, 8).Value = 263.07
$ws.Cells.Item(270, 4).Value = 0
$ws.Cells.Item(270, 5).Value = 452.3
$ws.Cells.Item(272, 4).Value = 0
$ws.Cells.Item(272, 5).Value = 149.61000000000001
$ws.Cells.Item(272, 8).Value = 149.61000000000001
$ws.Cells.Item(274, 4).Value = 0
$ws.Cells.Item(274, 5).Value = 383.23
$ws.Cells.Item(274, 8).Value = 383.23
